# Change the table style on the cash-flow table (slide 16) from the
# custom "Table_0" style to the built-in table style
# {7764F5CC-293D-4AA5-BFD4-EAC2D3582482}.
#
# Table styles can't be assigned through the Table.Style property
# (PowerPoint surfaces that as a read-only value) - the engine itself
# tells you to use Table.ApplyStyle(styleId) instead, which is the same
# call PowerPoint's Table Design ribbon issues when you click a style
# swatch in the gallery.

$p = $ppt.ActivePresentation

$slide = $p.Slides.Item(16)

$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $tableShape = $shape
        break
    }
}

$table = $tableShape.Table
$table.ApplyStyle("{7764F5CC-293D-4AA5-BFD4-EAC2D3582482}")
